$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9 (shifts "Aim of module" and everything below down by one)
$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = "Not Taken"
$ws.Range("B9").Value = "TP 500 Team Project Full"

$ws.Range("B9").Select()
